# Updates the "想去人数" (interested-attendee count) column F values
# on sheets 展览, 演出, and 全部类型 to match the regenerated data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7707
$ws.Range("F3").Value = 7707
$ws.Range("F5").Value = 7883
$ws.Range("F6").Value = 41
$ws.Range("F9").Value = 6663
$ws.Range("F10").Value = 3385
$ws.Range("F14").Value = 48
$ws.Range("F20").Value = 37
$ws.Range("F24").Value = 3855
$ws.Range("F28").Value = 283
$ws.Range("F29").Value = 1485
$ws.Range("F31").Value = 57
$ws.Range("F32").Value = 2760
$ws.Range("F33").Value = 1837
$ws.Range("F36").Value = 60
$ws.Range("F37").Value = 3692
$ws.Range("F38").Value = 320
$ws.Range("F41").Value = 919
$ws.Range("F43").Value = 1
$ws.Range("F47").Value = 551
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 409
$ws.Range("F7").Value = 42
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 7707
$ws.Range("F6").Value = 7707
$ws.Range("F7").Value = 7883
$ws.Range("F8").Value = 41
$ws.Range("F10").Value = 6663
$ws.Range("F11").Value = 3385
$ws.Range("F13").Value = 48
$ws.Range("F18").Value = 42
$ws.Range("F19").Value = 37
$ws.Range("F23").Value = 3855
$ws.Range("F29").Value = 283
$ws.Range("F30").Value = 1485
$ws.Range("F32").Value = 57
$ws.Range("F33").Value = 2760
$ws.Range("F34").Value = 1837
$ws.Range("F37").Value = 60
$ws.Range("F38").Value = 3692
$ws.Range("F39").Value = 320
$ws.Range("F43").Value = 919
$ws.Range("F49").Value = 551
